$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CreatedAt timestamp in A1
$ws.Range("A1").Value = "CreatedAt: 2025-12-06T19:06:53"

# Update numeric data cells (columns W, X, Y, Z) with refreshed IESO report values
$ws.Range("W4").Value = 131.46
$ws.Range("X4").Value = 163.61
$ws.Range("Y4").Value = 147.79
$ws.Range("Z4").Value = 147.8
$ws.Range("W6").Value = -5.92
$ws.Range("X6").Value = -6.71
$ws.Range("Y6").Value = -5.17
$ws.Range("Z6").Value = -3.7
$ws.Range("W9").Value = 131.84
$ws.Range("X9").Value = 166.81
$ws.Range("Y9").Value = 153.27
$ws.Range("Z9").Value = 152.72
$ws.Range("W11").Value = -5.54
$ws.Range("X11").Value = -3.5
$ws.Range("Y11").Value = 0.31
$ws.Range("Z11").Value = 1.22
$ws.Range("W14").Value = 131.84
$ws.Range("X14").Value = 166.81
$ws.Range("Y14").Value = 153.42
$ws.Range("Z14").Value = 152.72
$ws.Range("W16").Value = -5.54
$ws.Range("X16").Value = -3.5
$ws.Range("Y16").Value = 0.46
$ws.Range("Z16").Value = 1.22
$ws.Range("W19").Value = 129.84
$ws.Range("X19").Value = 161.74
$ws.Range("Y19").Value = 146.38
$ws.Range("Z19").Value = 146.52
$ws.Range("W21").Value = -7.53
$ws.Range("X21").Value = -8.57
$ws.Range("Y21").Value = -6.59
$ws.Range("Z21").Value = -4.98
$ws.Range("W24").Value = 129.84
$ws.Range("X24").Value = 161.74
$ws.Range("Y24").Value = 146.38
$ws.Range("Z24").Value = 146.52
$ws.Range("W26").Value = -7.53
$ws.Range("X26").Value = -8.57
$ws.Range("Y26").Value = -6.59
$ws.Range("Z26").Value = -4.98
$ws.Range("W29").Value = 127.67
$ws.Range("X29").Value = 159.32
$ws.Range("Y29").Value = 144.31
$ws.Range("Z29").Value = 144.42
$ws.Range("W31").Value = -9.699999999999999
$ws.Range("X31").Value = -10.99
$ws.Range("Y31").Value = -8.66
$ws.Range("Z31").Value = -7.08
$ws.Range("W34").Value = 134.95
$ws.Range("X34").Value = 173.26
$ws.Range("Y34").Value = 160.17
$ws.Range("Z34").Value = 158.97
$ws.Range("W36").Value = -2.43
$ws.Range("X36").Value = 2.95
$ws.Range("Y36").Value = 7.21
$ws.Range("Z36").Value = 7.47
$ws.Range("W39").Value = 131.46
$ws.Range("X39").Value = 163.61
$ws.Range("Y39").Value = 147.79
$ws.Range("Z39").Value = 147.8
$ws.Range("W41").Value = -5.92
$ws.Range("X41").Value = -6.71
$ws.Range("Y41").Value = -5.17
$ws.Range("Z41").Value = -3.7
$ws.Range("W44").Value = 138.2
$ws.Range("X44").Value = 171.35
$ws.Range("Y44").Value = 153.89
$ws.Range("Z44").Value = 153.34
$ws.Range("W46").Value = 0.83
$ws.Range("X46").Value = 1.03
$ws.Range("Y46").Value = 0.92
$ws.Range("Z46").Value = 1.84
$ws.Range("W49").Value = 146.93
$ws.Range("X49").Value = 181.19
$ws.Range("Y49").Value = 164.83
$ws.Range("Z49").Value = 162.21
$ws.Range("W51").Value = 9.550000000000001
$ws.Range("X51").Value = 10.87
$ws.Range("Y51").Value = 11.87
$ws.Range("Z51").Value = 10.71
$ws.Range("W54").Value = 137.1
$ws.Range("X54").Value = 171.69
$ws.Range("Y54").Value = 157.05
$ws.Range("Z54").Value = 156.83
$ws.Range("W56").Value = -0.27
$ws.Range("X56").Value = 1.37
$ws.Range("Y56").Value = 4.08
$ws.Range("Z56").Value = 5.33
$ws.Range("W59").Value = 143.1
$ws.Range("X59").Value = 177.04
$ws.Range("Y59").Value = 159.01
$ws.Range("Z59").Value = 157.98
$ws.Range("W61").Value = 5.72
$ws.Range("X61").Value = 6.73
$ws.Range("Y61").Value = 6.04
$ws.Range("Z61").Value = 6.48
$ws.Range("W64").Value = 145.68
$ws.Range("X64").Value = 180.23
$ws.Range("Y64").Value = 161.7
$ws.Range("Z64").Value = 160.66
$ws.Range("W66").Value = 8.300000000000001
$ws.Range("X66").Value = 9.91
$ws.Range("Y66").Value = 8.73
$ws.Range("Z66").Value = 9.16
$ws.Range("X69").Value = 777
$ws.Range("Y69").Value = 163.74
$ws.Range("Z69").Value = 162.54
$ws.Range("W71").Value = 9.710000000000001
$ws.Range("X71").Value = 11.84
$ws.Range("Y71").Value = 10.46
$ws.Range("Z71").Value = 10.88
$ws.Range("W72").Value = 357.92
$ws.Range("X72").Value = 594.84
$ws.Range("Y72").Value = 0.32
$ws.Range("Z72").Value = 0.16
$ws.Range("W74").Value = 144.45
$ws.Range("X74").Value = 178.9
$ws.Range("Y74").Value = 160.51
$ws.Range("Z74").Value = 159.31
$ws.Range("W76").Value = 7.08
$ws.Range("X76").Value = 8.59
$ws.Range("Y76").Value = 7.54
$ws.Range("Z76").Value = 7.81
$ws.Range("W79").Value = 145.05
$ws.Range("X79").Value = 179.48
$ws.Range("Y79").Value = 161
$ws.Range("Z79").Value = 159.84
$ws.Range("W81").Value = 7.67
$ws.Range("X81").Value = 9.16
$ws.Range("Y81").Value = 8.039999999999999
$ws.Range("Z81").Value = 8.34
$ws.Range("W84").Value = 131.84
$ws.Range("X84").Value = 165.2
$ws.Range("Y84").Value = 155.45
$ws.Range("Z84").Value = 154.12
$ws.Range("W86").Value = -5.54
$ws.Range("X86").Value = -5.12
$ws.Range("Y86").Value = 2.49
$ws.Range("Z86").Value = 2.62
$ws.Range("W89").Value = 127.67
$ws.Range("X89").Value = 159.32
$ws.Range("Y89").Value = 144.31
$ws.Range("Z89").Value = 144.42
$ws.Range("W91").Value = -9.699999999999999
$ws.Range("X91").Value = -10.99
$ws.Range("Y91").Value = -8.66
$ws.Range("Z91").Value = -7.08
